$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step A: trim "CoursesController (via Codegenerator)" down to just
# "CoursesController".
# ---------------------------------------------------------------------
$coursesPara = $d.Paragraphs.Item(72)
$rng = $coursesPara.Range
$rng.Find.Execute(" (via Codegenerator)", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------
# Step B: insert "LecturersController" and "StudentsController" list
# items right after the CoursesController item (continuing the same
# bullet list / numId).
# ---------------------------------------------------------------------
$coursesPara = $d.Paragraphs.Item(72)
$coursesPara.Range.InsertParagraphAfter()
$lectPara = $d.Paragraphs.Item(73)
$lr = $lectPara.Range
$lr.Collapse(1)
$lr.InsertAfter("Lecturers")
$lr.Collapse(0)
$lr.InsertAfter("Controller")

$lectPara2 = $d.Paragraphs.Item(73)
$lectPara2.Range.InsertParagraphAfter()
$studPara = $d.Paragraphs.Item(74)
$sr = $studPara.Range
$sr.Collapse(1)
$sr.InsertAfter("StudentsController")

# ---------------------------------------------------------------------
# Step C: remove the now-redundant blank paragraph that used to follow
# the CoursesController item.
# ---------------------------------------------------------------------
$stale = $d.Paragraphs.Item(75)
$stale.Range.Delete()

# ---------------------------------------------------------------------
# Step D: insert the new intro paragraph just before the
# "CoursesController" list item.
# ---------------------------------------------------------------------
$coursesPara = $d.Paragraphs.Item(72)
$coursesPara.Range.InsertParagraphBefore()

$introPara = $d.Paragraphs.Item(72)
$introPara.Style = "Normal"
$ir = $introPara.Range
$ir.Collapse(1)
$ir.InsertAfter("Created Controllers vis Codegenerator. ")
$ir.Collapse(0)
$ir.InsertAfter("Code Generator creates Controller Class, Razor view page")
$ir.Collapse(0)
$ir.InsertAfter(" (")
$ir.Collapse(0)
$ir.InsertAfter("Create, Edit, Details, Delete and Index page)")

# ---------------------------------------------------------------------
# Step E: after the CoursesController terminal command, add a blank
# Terminal paragraph, the LecturersController terminal command, another
# blank Terminal paragraph, the StudentsController terminal command,
# and a final blank Terminal paragraph.
# ---------------------------------------------------------------------
$cmdPara = $d.Paragraphs.Item(77)

$cmdPara.Range.InsertParagraphAfter()
$blank1 = $d.Paragraphs.Item(78)

$blank1.Range.InsertParagraphAfter()
$lectCmd = $d.Paragraphs.Item(79)
$lr2 = $lectCmd.Range
$lr2.Collapse(1)
$lr2.InsertAfter("dotnet aspnet-codegenerator controller -name ")
$lr2.Collapse(0)
$lr2.InsertAfter("Lecturers")
$lr2.Collapse(0)
$lr2.InsertAfter("Controller -m ")
$lr2.Collapse(0)
$lr2.InsertAfter("Lecturer")
$lr2.Collapse(0)
$lr2.InsertAfter(" -dc SchoolDbContext --relativeFolderPath Controllers --useDefaultLayout --referenceScriptLibraries -f")

$lectCmd2 = $d.Paragraphs.Item(79)
$lectCmd2.Range.InsertParagraphAfter()
$blank2 = $d.Paragraphs.Item(80)

$blank2.Range.InsertParagraphAfter()
$studCmd = $d.Paragraphs.Item(81)
$sr2 = $studCmd.Range
$sr2.Collapse(1)
$sr2.InsertAfter("dotnet aspnet-codegenerator controller -name ")
$sr2.Collapse(0)
$sr2.InsertAfter("Students")
$sr2.Collapse(0)
$sr2.InsertAfter("Controller -m ")
$sr2.Collapse(0)
$sr2.InsertAfter("Student")
$sr2.Collapse(0)
$sr2.InsertAfter(" -dc SchoolDbContext --relativeFolderPath Controllers --useDefaultLayout --referenceScriptLibraries -f")

$studCmd2 = $d.Paragraphs.Item(81)
$studCmd2.Range.InsertParagraphAfter()

Write-Host "Done"
